# Auto-generated edit script: updates crypto price/volume figures
# per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.329.08'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').Value = '3.325.04'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = "'584.23"
$ws.Range('E5').Value = '  +3.73%  '
$ws.Range('D6').Value = "'183.49"
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D9').Value = '3.321.01'
$ws.Range('E9').Value = '  +0.21%  '
$ws.Range('E10').Value = '  +1.78%  '
$ws.Range('D11').Value = "'0.582"
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').Value = "'46.45"
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = "'0.0000277"
$ws.Range('E13').Value = '  +4.74%  '
$ws.Range('D14').Value = "'645.14"
$ws.Range('E14').Value = '  +8.51%  '
$ws.Range('D15').Value = '3.856.30'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').Value = '68.386.21'
$ws.Range('E17').Value = '  +3.44%  '
$ws.Range('E18').Value = '  +1.59%  '
$ws.Range('D19').Value = '3.324.11'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = "'17.75"
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = "'10.95"
$ws.Range('D22').Value = "'0.904"
$ws.Range('E22').Value = '  +0.88%  '
$ws.Range('D23').Value = "'17.72"
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').Value = "'97.53"
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').Value = "'4.01"
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').Value = "'2.79"
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('D28').Value = "'9.62"
$ws.Range('E28').Value = '  +2.31%  '
$ws.Range('D29').Value = "'32.60"
$ws.Range('E29').Value = '  +6.32%  '
$ws.Range('D30').Value = "'8.62"
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('D31').Value = "'6.72"
$ws.Range('E31').Value = '  +1.21%  '
$ws.Range('D32').Value = "'603.15"
$ws.Range('E32').Value = '  +7.82%  '
$ws.Range('D33').Value = '3.965.09'
$ws.Range('E33').Value = '  +4.25%  '
$ws.Range('D34').Value = "'10.99"
$ws.Range('E34').Value = '  +1.25%  '
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').Value = "'3.53"
$ws.Range('E36').Value = '  -3.49%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('D38').Value = "'55.89"
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('D42').Value = "'32.71"
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').Value = '0.0₃0690'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('E44').Value = '  -0.23%  '
$ws.Range('D45').Value = "'0.338"
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('D46').Value = "'0.0417"
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').Value = "'2.56"
$ws.Range('E50').Value = '  +1.48%  '
$ws.Range('D51').Value = "'131.45"
$ws.Range('E51').Value = '  +1.91%  '

# Row 39/40 swap: Kaspa moves to row 39, Stacks moves to row 40
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = "'0.130"
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = "'3.28"
$ws.Range('E40').Value = '  +4.13%  '

# Row 47/48 swap: Stellar moves to row 47, Mantle moves to row 48
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = "'0.129"
$ws.Range('E47').Value = '  +1.91%  '
$ws.Range('B48').Value = 'Mantle'
$ws.Range('C48').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D48').Value = "'1.41"
$ws.Range('E48').Value = '  +13.88%  '
